# Automatische test-sync: 2025-08-30 19:03:50
#
# Adds the new "Klacht levering" log entry (row 5) to the Logs sheet,
# extends its conditional formatting ranges to include that row,
# adds the matching aggregate row (row 4) on the Dashboard sheet,
# and extends the chart's category/value series references accordingly.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row ---------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A5").Value = "Klacht levering"
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("D5").Value = "Klacht / Probleem"
$logs.Range("F5").Value = "2025-08-30 19:03:01"
$logs.Range("G5").Value = "Ja"
$logs.Range("H5").Value = "Nee"
$logs.Range("I5").Value = "Nee"
$logs.Range("J5").Value = "Nee"

# Extend the existing conditional formatting ranges (D, G, H, I, J)
# from row 2:4 to row 2:5, keeping their rules/dxf styles intact.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "4")
    $newRange = $logs.Range($col + "2:" + $col + "5")
    $fcs = $oldRange.FormatConditions
    $cnt = $fcs.Count()
    for ($i = 1; $i -le $cnt; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: append the aggregate row -----------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Klacht / Probleem"
$dash.Range("B4").Value = 1

# --- Chart: extend category/value series references -----------------
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.XValues = "=Dashboard!`$A`$2:`$A`$4"
$series.Values = "=Dashboard!`$B`$2:`$B`$4"
